$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SinhVien")

# Row 3: fill in Gender (G) and Marital/NU status (H) columns
$ws.Range("G3").Value = "NO"
$ws.Range("H3").Value = "NU"

# Row 4: fill in G column, add new H value "Nu"
$ws.Range("G4").Value = "NO"
$ws.Range("H4").Value = "Nu"

# Row 5: add G and H values
$ws.Range("G5").Value = "NO"
$ws.Range("H5").Value = "NU"

# Row 6: brand new row of data
$ws.Range("A6").Value = 113
$ws.Range("B6").Value = "ddsadas"
$ws.Range("C6").Value = 123
$ws.Range("D6").Value = "dsa"
$ws.Range("F6").Value = 44120.787280092591
$ws.Range("H6").Value = "NAM"
$ws.Range("I6").Value = "dsadas"

$ws.Range("G6").Select() | Out-Null
